# Adding Search Field component
# Append a new row (row 15) describing a "SearchField" component, following
# the same key/value layout used by the other component rows in the sheet.
#
# Values are written in the specific order below so that newly introduced
# shared strings land in the workbook's shared-strings table in the same
# order as the source edit (SearchField, sicCode, SIC Code, SIC Code:,
# placeHolder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value2 = "SearchField"
$ws.Range("C15").Value2 = "id"
$ws.Range("D15").Value2 = "sicCode"
$ws.Range("E15").Value2 = "label"
$ws.Range("J15").Value2 = "SIC Code"
$ws.Range("F15").Value2 = "SIC Code:"
$ws.Range("G15").Value2 = "mandatory"
$ws.Range("H15").Value2 = $false
$ws.Range("I15").Value2 = "placeHolder"

# Leave the selection where the author ended up after entering the row.
[void]$ws.Range("I16").Select()
